$d = $word.ActiveDocument

function Replace-ParagraphXml {
    param(
        [string]$AnchorText,
        [string]$ParagraphXml
    )
    $rng = $d.Content
    $rng.Find.Execute($AnchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $rng.Find.Found) {
        throw "Anchor text not found: $AnchorText"
    }
    $para = $rng.Paragraphs(1).Range
    $pkg = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' + $ParagraphXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $para.InsertXML($pkg)
}

# 1) "Manuel Dias" -> "Manuel Inácio Veladas Dias"
$d.Content.Find.Execute("Manuel Dias", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Manuel Inácio Veladas Dias", 2)

# 2) Remove proofErr spell-check wrapper around "help-desk" and merge the
#    surrounding runs (same rPr) into a single run.
$para3 = '<w:p w14:paraId="5AFA54BB" w14:textId="77777777" w:rsidR="00616450" w:rsidRPr="007B4B7B" w:rsidRDefault="00616450" w:rsidP="00616450"><w:pPr><w:widowControl w:val="0"/><w:spacing w:after="120" w:line="360" w:lineRule="auto"/><w:ind w:left="284" w:hanging="284"/><w:jc w:val="both"/></w:pPr><w:r w:rsidRPr="007B4B7B"><w:t>b)</w:t></w:r><w:r w:rsidRPr="007B4B7B"><w:tab/><w:t>Garantir a administração, operação, help-desk e manutenção do Fornecedor de Autenticação (Autenticação.Gov) e serviços de assinatura da Chave Móvel Digital;</w:t></w:r></w:p>'
Replace-ParagraphXml -AnchorText "help-desk" -ParagraphXml $para3

# 3) Remove proofErr spell-check wrapper around the first "guidelines"
#    occurrence (run itself is left untouched - different rPr from neighbours).
$para1 = '<w:p w14:paraId="74DC923C" w14:textId="64AA81BA" w:rsidR="00F901D1" w:rsidRPr="007B4B7B" w:rsidRDefault="00B36A41" w:rsidP="00F901D1"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="120" w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r w:rsidRPr="007B4B7B"><w:rPr><w:rFonts w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">No âmbito da assinatura com Chave Móvel Digital, cumprir as </w:t></w:r><w:r w:rsidRPr="007B4B7B"><w:rPr><w:rFonts w:cs="Times New Roman"/><w:i/></w:rPr><w:t>guidelines</w:t></w:r><w:r w:rsidRPr="007B4B7B"><w:rPr><w:rFonts w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> para implementação de aplicação de assinatura constantes da documentação disponibilizada pela </w:t></w:r><w:r w:rsidR="00407D79"><w:rPr><w:rFonts w:cs="Times New Roman"/></w:rPr><w:t>ARTE</w:t></w:r><w:r w:rsidRPr="007B4B7B"><w:rPr><w:rFonts w:cs="Times New Roman"/></w:rPr><w:t>;</w:t></w:r></w:p>'
Replace-ParagraphXml -AnchorText "cumprir as" -ParagraphXml $para1

# 4) Remove proofErr spell-check wrapper around the second "guidelines"
#    occurrence and merge "guidelines" with the following space run (both
#    share the same italic rPr).
$para2 = '<w:p w14:paraId="4CE30CCB" w14:textId="2F4EC23F" w:rsidR="00F901D1" w:rsidRPr="007B4B7B" w:rsidRDefault="00B36A41" w:rsidP="00F901D1"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="120" w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r w:rsidRPr="007B4B7B"><w:rPr><w:rFonts w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">No âmbito da assinatura com Chave Móvel Digital, disponibilizar à </w:t></w:r><w:r w:rsidR="00407D79"><w:rPr><w:rFonts w:cs="Times New Roman"/></w:rPr><w:t>ARTE</w:t></w:r><w:r w:rsidRPr="007B4B7B"><w:rPr><w:rFonts w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> documento que demonstre, para cada uma das </w:t></w:r><w:r w:rsidRPr="007B4B7B"><w:rPr><w:rFonts w:cs="Times New Roman"/><w:i/></w:rPr><w:t xml:space="preserve">guidelines </w:t></w:r><w:r w:rsidRPr="007B4B7B"><w:rPr><w:rFonts w:cs="Times New Roman"/></w:rPr><w:t>definidas, evidências do seu cumprimento;</w:t></w:r></w:p>'
Replace-ParagraphXml -AnchorText "cada uma das" -ParagraphXml $para2

# 5) Style DefaultParagraphFont becomes semiHidden (matches Word's default
#    "Default Paragraph Font" definition when the style sheet is refreshed).
$dpf = $d.Styles("Default Paragraph Font")
$dpf.SemiHidden = $true
